$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.015.19"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.047.22"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'246.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'0.658"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'57.24"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'0.0773"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "'15.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "'0.895"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.02%  "
$ws.Range("D14").Value = "2.349.34"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "'5.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "2.057.04"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'18.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.81%  "
$ws.Range("D18").Value = "36.985.41"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "'74.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'5.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "'236.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.10%  "
$ws.Range("D25").Value = "'9.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.72%  "
$ws.Range("D26").Value = "'169.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").Value = "'2.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "'19.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'5.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.24%  "
$ws.Range("D30").Value = "'0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "'4.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.14%  "
$ws.Range("D33").Value = "'0.0616"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("B34").Value = "BinanceUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.0874"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "'2.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("D37").Value = "'1.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.15%  "
$ws.Range("D38").Value = "'1.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").Value = "'5.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("D41").Value = "'0.0993"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.29%  "
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").Value = "'1.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("D44").Value = "'97.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").Value = "'17.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.293.94"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'2.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "'6.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  +6.47%  "
$ws.Range("D51").Value = "2.234.97"
$ws.Range("E51").Value = "  +0.01%  "
